# Scheduled-runner data refresh: updates computed market/profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) on the rows
# whose underlying item prices moved, across the ALC/ARM/BSM/CRP/CUL/GSM/
# LTW/WVR leve-profit sheets. Row/column layout and item data are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1803.3077  # H19: 2145.5 -> 1803.3077
$ws.Cells.Item(19, 9).Value = 841.1111  # I19: 856.8333 -> 841.1111
$ws.Cells.Item(19, 10).Value = 3968.25  # J19: 4078.5 -> 3968.25
$ws.Cells.Item(19, 11).Value = 841.1111  # K19: 856.8333 -> 841.1111
$ws.Cells.Item(19, 12).Value = 3968.25  # L19: 4078.5 -> 3968.25
$ws.Cells.Item(19, 13).Value = -666.1111  # M19: -681.8333 -> -666.1111
$ws.Cells.Item(19, 14).Value = -4318.25  # N19: -4428.5 -> -4318.25

$ws.Cells.Item(32, 8).Value = 5031.8335  # H32: 4299 -> 5031.8335
$ws.Cells.Item(32, 10).Value = 6110.5557  # J32: 5371.857 -> 6110.5557
$ws.Cells.Item(32, 12).Value = 6110.5557  # L32: 5371.857 -> 6110.5557
$ws.Cells.Item(32, 14).Value = -6762.5557  # N32: -6023.857 -> -6762.5557

$ws.Cells.Item(33, 8).Value = 256.5  # H33: 390.82352 -> 256.5
$ws.Cells.Item(33, 9).Value = 272.85715  # I33: 275.7857 -> 272.85715
$ws.Cells.Item(33, 10).Value = 142  # J33: 927.6667 -> 142
$ws.Cells.Item(33, 11).Value = 272.85715  # K33: 275.7857 -> 272.85715
$ws.Cells.Item(33, 12).Value = 142  # L33: 927.6667 -> 142
$ws.Cells.Item(33, 13).Value = -43.85714999999999  # M33: -46.78570000000002 -> -43.85714999999999
$ws.Cells.Item(33, 14).Value = -600  # N33: -1385.6667 -> -600

$ws.Cells.Item(40, 8).Value = 7146151.5  # H40: 13574548 -> 7146151.5
$ws.Cells.Item(40, 9).Value = 1987604.5  # I40: 3971720.2 -> 1987604.5
$ws.Cells.Item(40, 11).Value = 1987604.5  # K40: 3971720.2 -> 1987604.5
$ws.Cells.Item(40, 13).Value = -1987429.5  # M40: -3971545.2 -> -1987429.5

$ws.Cells.Item(64, 8).Value = 7422.846  # H64: 7541.5835 -> 7422.846
$ws.Cells.Item(64, 10).Value = 7599.4  # J64: 7999.75 -> 7599.4
$ws.Cells.Item(64, 12).Value = 7599.4  # L64: 7999.75 -> 7599.4
$ws.Cells.Item(64, 14).Value = -8095.4  # N64: -8495.75 -> -8095.4

$ws.Cells.Item(67, 8).Value = 7422.846  # H67: 7541.5835 -> 7422.846
$ws.Cells.Item(67, 10).Value = 7599.4  # J67: 7999.75 -> 7599.4
$ws.Cells.Item(67, 12).Value = 7599.4  # L67: 7999.75 -> 7599.4
$ws.Cells.Item(67, 14).Value = -9315.4  # N67: -9715.75 -> -9315.4

$ws.Cells.Item(113, 8).Value = 11062.125  # H113: 11999.857 -> 11062.125
$ws.Cells.Item(113, 9).Value = 7124.25  # I113: 7999.6665 -> 7124.25
$ws.Cells.Item(113, 11).Value = 7124.25  # K113: 7999.6665 -> 7124.25
$ws.Cells.Item(113, 13).Value = -3870.25  # M113: -4745.6665 -> -3870.25

$ws.Cells.Item(132, 8).Value = 2631.3704  # H132: 2635 -> 2631.3704
$ws.Cells.Item(132, 9).Value = 2658.5833  # I132: 2651.96 -> 2658.5833
$ws.Cells.Item(132, 10).Value = 2413.6667  # J132: 2423 -> 2413.6667
$ws.Cells.Item(132, 11).Value = 7975.749899999999  # K132: 7955.88 -> 7975.749899999999
$ws.Cells.Item(132, 12).Value = 7241.000100000001  # L132: 7269 -> 7241.000100000001
$ws.Cells.Item(132, 13).Value = -5445.749899999999  # M132: -5425.88 -> -5445.749899999999
$ws.Cells.Item(132, 14).Value = -12301.0001  # N132: -12329 -> -12301.0001

$ws.Cells.Item(138, 8).Value = 2976.5208  # H138: 3452.3489 -> 2976.5208
$ws.Cells.Item(138, 9).Value = 1584.7273  # I138: 1654 -> 1584.7273
$ws.Cells.Item(138, 10).Value = 3390.2974  # J138: 3997.303 -> 3390.2974
$ws.Cells.Item(138, 11).Value = 4754.1819  # K138: 4962 -> 4754.1819
$ws.Cells.Item(138, 12).Value = 10170.8922  # L138: 11991.909 -> 10170.8922
$ws.Cells.Item(138, 13).Value = 385.8181000000004  # M138: 178 -> 385.8181000000004
$ws.Cells.Item(138, 14).Value = -20450.8922  # N138: -22271.909 -> -20450.8922

$ws.Cells.Item(141, 8).Value = 8088.8823  # H141: 8123.5884 -> 8088.8823
$ws.Cells.Item(141, 9).Value = 6744.5454  # I141: 6798.1816 -> 6744.5454
$ws.Cells.Item(141, 11).Value = 20233.6362  # K141: 20394.5448 -> 20233.6362
$ws.Cells.Item(141, 13).Value = -15053.6362  # M141: -15214.5448 -> -15053.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2868.375  # H61: 2927.7144 -> 2868.375
$ws.Cells.Item(61, 9).Value = 1507.4445  # I61: 1521.875 -> 1507.4445
$ws.Cells.Item(61, 10).Value = 4618.143  # J61: 4802.1665 -> 4618.143
$ws.Cells.Item(61, 11).Value = 1507.4445  # K61: 1521.875 -> 1507.4445
$ws.Cells.Item(61, 12).Value = 4618.143  # L61: 4802.1665 -> 4618.143
$ws.Cells.Item(61, 13).Value = -1295.4445  # M61: -1309.875 -> -1295.4445
$ws.Cells.Item(61, 14).Value = -5042.143  # N61: -5226.1665 -> -5042.143

$ws.Cells.Item(136, 8).Value = 2868.375  # H136: 2927.7144 -> 2868.375
$ws.Cells.Item(136, 9).Value = 1507.4445  # I136: 1521.875 -> 1507.4445
$ws.Cells.Item(136, 10).Value = 4618.143  # J136: 4802.1665 -> 4618.143
$ws.Cells.Item(136, 11).Value = 4522.333500000001  # K136: 4565.625 -> 4522.333500000001
$ws.Cells.Item(136, 12).Value = 13854.429  # L136: 14406.4995 -> 13854.429
$ws.Cells.Item(136, 13).Value = -1972.333500000001  # M136: -2015.625 -> -1972.333500000001
$ws.Cells.Item(136, 14).Value = -18954.429  # N136: -19506.4995 -> -18954.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(110, 8).Value = 90000  # H110: 39900 -> 90000
$ws.Cells.Item(110, 10).Value = 90000  # J110: 39900 -> 90000
$ws.Cells.Item(110, 12).Value = 90000  # L110: 39900 -> 90000
$ws.Cells.Item(110, 14).Value = -98180  # N110: -48080 -> -98180

$ws.Cells.Item(134, 8).Value = 3238.4736  # H134: 3508.3125 -> 3238.4736
$ws.Cells.Item(134, 9).Value = 3238.4736  # I134: 3508.3125 -> 3238.4736
$ws.Cells.Item(134, 11).Value = 9715.4208  # K134: 10524.9375 -> 9715.4208
$ws.Cells.Item(134, 13).Value = -7180.4208  # M134: -7989.9375 -> -7180.4208

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11, 8).Value = 1348.5  # H11: 1532.3334 -> 1348.5
$ws.Cells.Item(11, 10).Value = 1348.5  # J11: 1532.3334 -> 1348.5
$ws.Cells.Item(11, 12).Value = 1348.5  # L11: 1532.3334 -> 1348.5
$ws.Cells.Item(11, 14).Value = -1628.5  # N11: -1812.3334 -> -1628.5

$ws.Cells.Item(99, 8).Value = 1843.6666  # H99: 1869.5135 -> 1843.6666
$ws.Cells.Item(99, 9).Value = 1811.6562  # I99: 1811.6875 -> 1811.6562
$ws.Cells.Item(99, 10).Value = 2099.75  # J99: 2239.6 -> 2099.75
$ws.Cells.Item(99, 11).Value = 1811.6562  # K99: 1811.6875 -> 1811.6562
$ws.Cells.Item(99, 12).Value = 2099.75  # L99: 2239.6 -> 2099.75
$ws.Cells.Item(99, 13).Value = -313.6561999999999  # M99: -313.6875 -> -313.6561999999999
$ws.Cells.Item(99, 14).Value = -5095.75  # N99: -5235.6 -> -5095.75

$ws.Cells.Item(105, 8).Value = 2740.111  # H105: 3031.5334 -> 2740.111
$ws.Cells.Item(105, 9).Value = 2718.1  # I105: 3106.5 -> 2718.1
$ws.Cells.Item(105, 10).Value = 2767.625  # J105: 2945.8572 -> 2767.625
$ws.Cells.Item(105, 11).Value = 2718.1  # K105: 3106.5 -> 2718.1
$ws.Cells.Item(105, 12).Value = 2767.625  # L105: 2945.8572 -> 2767.625
$ws.Cells.Item(105, 13).Value = -971.0999999999999  # M105: -1359.5 -> -971.0999999999999
$ws.Cells.Item(105, 14).Value = -6261.625  # N105: -6439.8572 -> -6261.625

$ws.Cells.Item(122, 8).Value = 1249.7142  # H122: 1335.4286 -> 1249.7142
$ws.Cells.Item(122, 9).Value = 800  # I122: 950 -> 800
$ws.Cells.Item(122, 10).Value = 1429.6  # J122: 1399.6666 -> 1429.6
$ws.Cells.Item(122, 11).Value = 2400  # K122: 2850 -> 2400
$ws.Cells.Item(122, 12).Value = 4288.799999999999  # L122: 4198.9998 -> 4288.799999999999
$ws.Cells.Item(122, 13).Value = 50  # M122: -400 -> 50
$ws.Cells.Item(122, 14).Value = -9188.799999999999  # N122: -9098.9998 -> -9188.799999999999

$ws.Cells.Item(126, 8).Value = 1843.6666  # H126: 1869.5135 -> 1843.6666
$ws.Cells.Item(126, 9).Value = 1811.6562  # I126: 1811.6875 -> 1811.6562
$ws.Cells.Item(126, 10).Value = 2099.75  # J126: 2239.6 -> 2099.75
$ws.Cells.Item(126, 11).Value = 5434.9686  # K126: 5435.0625 -> 5434.9686
$ws.Cells.Item(126, 12).Value = 6299.25  # L126: 6718.799999999999 -> 6299.25
$ws.Cells.Item(126, 13).Value = -2964.9686  # M126: -2965.0625 -> -2964.9686
$ws.Cells.Item(126, 14).Value = -11239.25  # N126: -11658.8 -> -11239.25

$ws.Cells.Item(134, 8).Value = 4434.3335  # H134: 4462.385 -> 4434.3335
$ws.Cells.Item(134, 9).Value = 4341.2  # I134: 4382.8184 -> 4341.2
$ws.Cells.Item(134, 11).Value = 13023.6  # K134: 13148.4552 -> 13023.6
$ws.Cells.Item(134, 13).Value = -10488.6  # M134: -10613.4552 -> -10488.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3369.75  # H39: 3485.6365 -> 3369.75
$ws.Cells.Item(39, 10).Value = 3779.7144  # J39: 3784.2 -> 3779.7144
$ws.Cells.Item(39, 12).Value = 11339.1432  # L39: 11352.6 -> 11339.1432
$ws.Cells.Item(39, 14).Value = -11927.1432  # N39: -11940.6 -> -11927.1432

$ws.Cells.Item(62, 8).Value = 4099  # H62: 3199 -> 4099
$ws.Cells.Item(62, 9).Value = 4999  # I62: 0 -> 4999
$ws.Cells.Item(62, 11).Value = 14997  # K62: 0 -> 14997
$ws.Cells.Item(62, 13).Value = -14311  # M62: None -> -14311

$ws.Cells.Item(65, 8).Value = 4099  # H65: 3199 -> 4099
$ws.Cells.Item(65, 9).Value = 4999  # I65: 0 -> 4999
$ws.Cells.Item(65, 11).Value = 44991  # K65: 0 -> 44991
$ws.Cells.Item(65, 13).Value = -41559  # M65: None -> -41559

$ws.Cells.Item(92, 8).Value = 324.2857  # H92: 318.46667 -> 324.2857
$ws.Cells.Item(92, 9).Value = 249.66667  # I92: 246.5 -> 249.66667
$ws.Cells.Item(92, 11).Value = 749.00001  # K92: 739.5 -> 749.00001
$ws.Cells.Item(92, 13).Value = 498.99999  # M92: 508.5 -> 498.99999

$ws.Cells.Item(138, 8).Value = 4002.3  # H138: 3660.9092 -> 4002.3
$ws.Cells.Item(138, 9).Value = 2780.3333  # I138: 2527 -> 2780.3333
$ws.Cells.Item(138, 11).Value = 8340.999899999999  # K138: 7581 -> 8340.999899999999
$ws.Cells.Item(138, 13).Value = -3200.999899999999  # M138: -2441 -> -3200.999899999999

$ws.Cells.Item(140, 8).Value = 3488.5217  # H140: 3633.4546 -> 3488.5217
$ws.Cells.Item(140, 9).Value = 1011.8  # I140: 1049.2632 -> 1011.8
$ws.Cells.Item(140, 11).Value = 3035.4  # K140: 3147.7896 -> 3035.4
$ws.Cells.Item(140, 13).Value = 2144.6  # M140: 2032.2104 -> 2144.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 0  # H104: 25000 -> 0
$ws.Cells.Item(104, 10).Value = 0  # J104: 25000 -> 0
$ws.Cells.Item(104, 12).Value = 0  # L104: 25000 -> 0
$ws.Cells.Item(104, 14).ClearContents()  # N104 remove (was -31988)

$ws.Cells.Item(122, 8).Value = 2317.6667  # H122: 2383.2917 -> 2317.6667
$ws.Cells.Item(122, 9).Value = 2270.1052  # I122: 2295.35 -> 2270.1052
$ws.Cells.Item(122, 10).Value = 2498.4  # J122: 2823 -> 2498.4
$ws.Cells.Item(122, 11).Value = 6810.3156  # K122: 6886.049999999999 -> 6810.3156
$ws.Cells.Item(122, 12).Value = 7495.200000000001  # L122: 8469 -> 7495.200000000001
$ws.Cells.Item(122, 13).Value = -4360.3156  # M122: -4436.049999999999 -> -4360.3156
$ws.Cells.Item(122, 14).Value = -12395.2  # N122: -13369 -> -12395.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 6946305  # H46: 6946324 -> 6946305
$ws.Cells.Item(46, 9).Value = 16669527  # I46: 20836808 -> 16669527
$ws.Cells.Item(46, 10).Value = 1146.1428  # J46: 1081.625 -> 1146.1428
$ws.Cells.Item(46, 11).Value = 16669527  # K46: 20836808 -> 16669527
$ws.Cells.Item(46, 12).Value = 1146.1428  # L46: 1081.625 -> 1146.1428
$ws.Cells.Item(46, 13).Value = -16669339  # M46: -20836620 -> -16669339
$ws.Cells.Item(46, 14).Value = -1522.1428  # N46: -1457.625 -> -1522.1428

$ws.Cells.Item(61, 8).Value = 5241.2173  # H61: 5143.625 -> 5241.2173
$ws.Cells.Item(61, 9).Value = 5392.857  # I61: 5279.5 -> 5392.857
$ws.Cells.Item(61, 11).Value = 5392.857  # K61: 5279.5 -> 5392.857
$ws.Cells.Item(61, 13).Value = -5190.857  # M61: -5077.5 -> -5190.857

$ws.Cells.Item(113, 8).Value = 5241.2173  # H113: 5143.625 -> 5241.2173
$ws.Cells.Item(113, 9).Value = 5392.857  # I113: 5279.5 -> 5392.857
$ws.Cells.Item(113, 11).Value = 5392.857  # K113: 5279.5 -> 5392.857
$ws.Cells.Item(113, 13).Value = -3222.857  # M113: -3109.5 -> -3222.857

$ws.Cells.Item(132, 8).Value = 5056.909  # H132: 4939.385 -> 5056.909
$ws.Cells.Item(132, 10).Value = 4398  # J132: 4328 -> 4398
$ws.Cells.Item(132, 12).Value = 13194  # L132: 12984 -> 13194
$ws.Cells.Item(132, 14).Value = -18254  # N132: -18044 -> -18254

$ws.Cells.Item(136, 8).Value = 3712.1765  # H136: 3713.7058 -> 3712.1765
$ws.Cells.Item(136, 9).Value = 3073.1333  # I136: 3074.8667 -> 3073.1333
$ws.Cells.Item(136, 11).Value = 9219.3999  # K136: 9224.6001 -> 9219.3999
$ws.Cells.Item(136, 13).Value = -6669.3999  # M136: -6674.6001 -> -6669.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 19234180  # H122: 20003506 -> 19234180
$ws.Cells.Item(122, 9).Value = 26319426  # I122: 27781560 -> 26319426
$ws.Cells.Item(122, 11).Value = 78958278  # K122: 83344680 -> 78958278
$ws.Cells.Item(122, 13).Value = -78955828  # M122: -83342230 -> -78955828

$ws.Cells.Item(132, 8).Value = 1536.9375  # H132: 1533.6086 -> 1536.9375
$ws.Cells.Item(132, 9).Value = 1439.4  # I132: 1363.85 -> 1439.4
$ws.Cells.Item(132, 10).Value = 3000  # J132: 2665.3333 -> 3000
$ws.Cells.Item(132, 11).Value = 4318.200000000001  # K132: 4091.55 -> 4318.200000000001
$ws.Cells.Item(132, 12).Value = 9000  # L132: 7995.999899999999 -> 9000
$ws.Cells.Item(132, 13).Value = -1788.200000000001  # M132: -1561.55 -> -1788.200000000001
$ws.Cells.Item(132, 14).Value = -14060  # N132: -13055.9999 -> -14060

$ws.Cells.Item(133, 8).Value = 92616.664  # H133: 86274.75 -> 92616.664
$ws.Cells.Item(133, 10).Value = 108600  # J133: 94816.336 -> 108600
$ws.Cells.Item(133, 12).Value = 108600  # L133: 94816.336 -> 108600
$ws.Cells.Item(133, 14).Value = -118720  # N133: -104936.336 -> -118720
